# Auto-generated edit script applying numeric updates from the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 599.8077
$ws.Range("I15").Value = 599.8077
$ws.Range("K15").Value = 1799.4231
$ws.Range("M15").Value = -1630.4231

$ws.Range("H38").Value = 854.8
$ws.Range("I38").Value = 854.8
$ws.Range("K38").Value = 2564.4
$ws.Range("M38").Value = -2192.4

$ws.Range("H62").Value = 125005830
$ws.Range("I62").Value = 142859800
$ws.Range("K62").Value = 142859800
$ws.Range("M62").Value = -142859176

$ws.Range("H65").Value = 125005830
$ws.Range("I65").Value = 142859800
$ws.Range("K65").Value = 714299000
$ws.Range("M65").Value = -714295880

$ws.Range("H74").Value = 3879.6
$ws.Range("I74").Value = 2199.5
$ws.Range("K74").Value = 2199.5
$ws.Range("M74").Value = -1263.5

$ws.Range("H76").Value = 3433
$ws.Range("J76").Value = 3433
$ws.Range("L76").Value = 3433
$ws.Range("N76").Value = -4063

$ws.Range("H77").Value = 3879.6
$ws.Range("I77").Value = 2199.5
$ws.Range("K77").Value = 10997.5
$ws.Range("M77").Value = -6317.5

$ws.Range("H79").Value = 3433
$ws.Range("J79").Value = 3433
$ws.Range("L79").Value = 3433
$ws.Range("N79").Value = -5617

$ws.Range("H132").Value = 864.27026
$ws.Range("I132").Value = 754.89655
$ws.Range("K132").Value = 2264.68965
$ws.Range("M132").Value = 265.3103499999997

$ws.Range("H141").Value = 3115523
$ws.Range("I141").Value = 4001783.5
$ws.Range("J141").Value = 13611
$ws.Range("K141").Value = 12005350.5
$ws.Range("L141").Value = 40833
$ws.Range("M141").Value = -12000170.5
$ws.Range("N141").Value = -51193

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 309930.12
$ws.Range("I2").Value = 464171.75
$ws.Range("K2").Value = 464171.75
$ws.Range("M2").Value = -464058.75

$ws.Range("H32").Value = 4293.7573
$ws.Range("I32").Value = 2787.5938
$ws.Range("K32").Value = 2787.5938
$ws.Range("M32").Value = -2500.5938

$ws.Range("H45").Value = 6430050.5
$ws.Range("I45").Value = 18000802
$ws.Range("K45").Value = 18000802
$ws.Range("M45").Value = -18000425

$ws.Range("H61").Value = 6416.8696
$ws.Range("I61").Value = 7126.1333
$ws.Range("J61").Value = 5087
$ws.Range("K61").Value = 7126.1333
$ws.Range("L61").Value = 5087
$ws.Range("M61").Value = -6914.1333
$ws.Range("N61").Value = -5511

$ws.Range("H63").Value = 1702.2
$ws.Range("I63").Value = 1652.75
$ws.Range("J63").Value = 1900
$ws.Range("K63").Value = 1652.75
$ws.Range("L63").Value = 1900
$ws.Range("M63").Value = -966.75
$ws.Range("N63").Value = -3272

$ws.Range("H66").Value = 1702.2
$ws.Range("I66").Value = 1652.75
$ws.Range("J66").Value = 1900
$ws.Range("K66").Value = 8263.75
$ws.Range("L66").Value = 9500
$ws.Range("M66").Value = -4831.75
$ws.Range("N66").Value = -16364

$ws.Range("H74").Value = 1378.9166
$ws.Range("I74").Value = 500
$ws.Range("J74").Value = 3136.75
$ws.Range("K74").Value = 500
$ws.Range("L74").Value = 3136.75
$ws.Range("M74").Value = 374
$ws.Range("N74").Value = -4884.75

$ws.Range("H77").Value = 1378.9166
$ws.Range("I77").Value = 500
$ws.Range("J77").Value = 3136.75
$ws.Range("K77").Value = 2500
$ws.Range("L77").Value = 15683.75
$ws.Range("M77").Value = 1868
$ws.Range("N77").Value = -24419.75

$ws.Range("H116").Value = 309930.12
$ws.Range("I116").Value = 464171.75
$ws.Range("K116").Value = 464171.75
$ws.Range("M116").Value = -461877.75

$ws.Range("H122").Value = 1288.5385
$ws.Range("I122").Value = 928
$ws.Range("J122").Value = 2099.75
$ws.Range("K122").Value = 2784
$ws.Range("L122").Value = 6299.25
$ws.Range("M122").Value = -334
$ws.Range("N122").Value = -11199.25

$ws.Range("H132").Value = 1750.1852
$ws.Range("I132").Value = 1312.25
$ws.Range("K132").Value = 3936.75
$ws.Range("M132").Value = -1406.75

$ws.Range("H136").Value = 6416.8696
$ws.Range("I136").Value = 7126.1333
$ws.Range("J136").Value = 5087
$ws.Range("K136").Value = 21378.3999
$ws.Range("L136").Value = 15261
$ws.Range("M136").Value = -18828.3999
$ws.Range("N136").Value = -20361

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 309930.12
$ws.Range("I3").Value = 464171.75
$ws.Range("K3").Value = 464171.75
$ws.Range("M3").Value = -464057.75

$ws.Range("H20").Value = 2586.7693
$ws.Range("I20").Value = 2392.3333
$ws.Range("K20").Value = 2392.3333
$ws.Range("M20").Value = -2145.3333

$ws.Range("H134").Value = 5244.517
$ws.Range("I134").Value = 5707
$ws.Range("K134").Value = 17121
$ws.Range("M134").Value = -14586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3333
$ws.Range("I86").Value = 3333
$ws.Range("K86").Value = 3333
$ws.Range("M86").Value = -2210

$ws.Range("H89").Value = 3333
$ws.Range("I89").Value = 3333
$ws.Range("K89").Value = 16665
$ws.Range("M89").Value = -11049

$ws.Range("H131").Value = 38886.215
$ws.Range("J131").Value = 38886.215
$ws.Range("L131").Value = 38886.215
$ws.Range("N131").Value = -48966.215

$ws.Range("H132").Value = 2424.2173
$ws.Range("I132").Value = 1543.2667
$ws.Range("J132").Value = 4076
$ws.Range("K132").Value = 4629.800099999999
$ws.Range("L132").Value = 12228
$ws.Range("M132").Value = -2099.800099999999
$ws.Range("N132").Value = -17288

$ws.Range("H134").Value = 1735.5714
$ws.Range("I134").Value = 1373.75
$ws.Range("K134").Value = 4121.25
$ws.Range("M134").Value = -1586.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 786
$ws.Range("J97").Value = 1075.75
$ws.Range("L97").Value = 3227.25
$ws.Range("N97").Value = -4219.25

$ws.Range("H122").Value = 1748.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1748.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 15736.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -20636.5

$ws.Range("H132").Value = 1540.875
$ws.Range("J132").Value = 2100.25
$ws.Range("L132").Value = 18902.25
$ws.Range("N132").Value = -23962.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 9000
$ws.Range("J55").Value = 9000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9654

$ws.Range("H97").Value = 570.74194
$ws.Range("I97").Value = 580.10345
$ws.Range("K97").Value = 580.10345
$ws.Range("M97").Value = -84.10344999999995

$ws.Range("H122").Value = 1841.5714
$ws.Range("I122").Value = 1763
$ws.Range("J122").Value = 1983
$ws.Range("K122").Value = 5289
$ws.Range("L122").Value = 5949
$ws.Range("M122").Value = -2839
$ws.Range("N122").Value = -10849

$ws.Range("H126").Value = 2830475
$ws.Range("I126").Value = 5558541
$ws.Range("K126").Value = 16675623
$ws.Range("M126").Value = -16673153

$ws.Range("H132").Value = 2026827.6
$ws.Range("I132").Value = 2748844.8
$ws.Range("K132").Value = 8246534.399999999
$ws.Range("M132").Value = -8244004.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2194.75
$ws.Range("J22").Value = 2693
$ws.Range("L22").Value = 2693
$ws.Range("N22").Value = -3283

$ws.Range("H27").Value = 2194.75
$ws.Range("J27").Value = 2693
$ws.Range("L27").Value = 2693
$ws.Range("N27").Value = -2907

$ws.Range("H61").Value = 2010.0454
$ws.Range("I61").Value = 1806.7142
$ws.Range("J61").Value = 2365.875
$ws.Range("K61").Value = 1806.7142
$ws.Range("L61").Value = 2365.875
$ws.Range("M61").Value = -1604.7142
$ws.Range("N61").Value = -2769.875

$ws.Range("H113").Value = 2010.0454
$ws.Range("I113").Value = 1806.7142
$ws.Range("J113").Value = 2365.875
$ws.Range("K113").Value = 1806.7142
$ws.Range("L113").Value = 2365.875
$ws.Range("M113").Value = 363.2858000000001
$ws.Range("N113").Value = -6705.875

$ws.Range("H132").Value = 1577.2046
$ws.Range("I132").Value = 1274.7142
$ws.Range("K132").Value = 3824.1426
$ws.Range("M132").Value = -1294.1426

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 33292.375
$ws.Range("I122").Value = 37748.523
$ws.Range("K122").Value = 113245.569
$ws.Range("M122").Value = -110795.569

$ws.Range("H132").Value = 1219.258
$ws.Range("I132").Value = 890.86957
$ws.Range("J132").Value = 2163.375
$ws.Range("K132").Value = 2672.60871
$ws.Range("L132").Value = 6490.125
$ws.Range("M132").Value = -142.60871
$ws.Range("N132").Value = -11550.125

$ws.Range("H136").Value = 1613.8182
$ws.Range("I136").Value = 1359
$ws.Range("J136").Value = 1919.6
$ws.Range("K136").Value = 4077
$ws.Range("L136").Value = 5758.799999999999
$ws.Range("M136").Value = -1527
$ws.Range("N136").Value = -10858.8
